$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H17").Value = 1805.76
$ws.Range("J17").Value = 1805.76
$ws.Range("L17").Value = 5417.28
$ws.Range("N17").Value = -5753.28

$ws.Range("H19").Value = 537.4
$ws.Range("I19").Value = 571.75
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 571.75
$ws.Range("L19").Value = 400
$ws.Range("M19").Value = -396.75
$ws.Range("N19").Value = -750

$ws.Range("H53").Value = 37.5
$ws.Range("I53").Value = 20.833334
$ws.Range("J53").Value = 62.5
$ws.Range("K53").Value = 20.833334
$ws.Range("L53").Value = 62.5
$ws.Range("M53").Value = 616.166666
$ws.Range("N53").Value = -1336.5

$ws.Range("H62").Value = 2735
$ws.Range("I62").Value = 2735
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2735
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2111
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 2735
$ws.Range("I65").Value = 2735
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 13675
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -10555
$ws.Range("N65").ClearContents()

$ws.Range("H112").Value = 2494.9092
$ws.Range("J112").Value = 2701.923
$ws.Range("L112").Value = 8105.768999999999
$ws.Range("N112").Value = -10321.769

$ws.Range("H137").Value = 1290.2131
$ws.Range("I137").Value = 1106.7192
$ws.Range("J137").Value = 3905
$ws.Range("K137").Value = 3320.1576
$ws.Range("L137").Value = 11715
$ws.Range("M137").Value = -770.1576
$ws.Range("N137").Value = -16815

$ws.Range("H138").Value = 2821722.8
$ws.Range("J138").Value = 7242.294
$ws.Range("L138").Value = 21726.882
$ws.Range("N138").Value = -32006.882

$ws.Range("H141").Value = 14324.375
$ws.Range("I141").Value = 2648.75
$ws.Range("J141").Value = 26000
$ws.Range("K141").Value = 7946.25
$ws.Range("L141").Value = 78000
$ws.Range("M141").Value = -2766.25
$ws.Range("N141").Value = -88360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 5625.636
$ws.Range("I31").Value = 5625.636
$ws.Range("K31").Value = 5625.636
$ws.Range("M31").Value = -5331.636

$ws.Range("H74").Value = 1287.0333
$ws.Range("I74").Value = 1314.3914
$ws.Range("J74").Value = 1197.1428
$ws.Range("K74").Value = 1314.3914
$ws.Range("L74").Value = 1197.1428
$ws.Range("M74").Value = -440.3914
$ws.Range("N74").Value = -2945.1428

$ws.Range("H77").Value = 1287.0333
$ws.Range("I77").Value = 1314.3914
$ws.Range("J77").Value = 1197.1428
$ws.Range("K77").Value = 6571.957
$ws.Range("L77").Value = 5985.714
$ws.Range("M77").Value = -2203.957
$ws.Range("N77").Value = -14721.714

$ws.Range("H102").Value = 2900
$ws.Range("I102").Value = 1516.6666
$ws.Range("J102").Value = 5666.6665
$ws.Range("K102").Value = 1516.6666
$ws.Range("L102").Value = 5666.6665
$ws.Range("M102").Value = 105.3334
$ws.Range("N102").Value = -8910.6665

$ws.Range("H132").Value = 911964.06
$ws.Range("I132").Value = 1252288.2
$ws.Range("K132").Value = 3756864.6
$ws.Range("M132").Value = -3754334.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 45593.332
$ws.Range("J60").Value = 79780
$ws.Range("L60").Value = 79780
$ws.Range("N60").Value = -80978

$ws.Range("H134").Value = 543388.7
$ws.Range("I134").Value = 716609.56
$ws.Range("J134").Value = 4479.3335
$ws.Range("K134").Value = 2149828.68
$ws.Range("L134").Value = 13438.0005
$ws.Range("M134").Value = -2147293.68
$ws.Range("N134").Value = -18508.0005

$ws.Range("H137").Value = 45936.363
$ws.Range("J137").Value = 45936.363
$ws.Range("L137").Value = 45936.363
$ws.Range("N137").Value = -56136.363

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2817.0986
$ws.Range("I31").Value = 1594.4147
$ws.Range("J31").Value = 4488.1
$ws.Range("K31").Value = 1594.4147
$ws.Range("L31").Value = 4488.1
$ws.Range("M31").Value = -1299.4147
$ws.Range("N31").Value = -5078.1

$ws.Range("H34").Value = 2817.0986
$ws.Range("I34").Value = 1594.4147
$ws.Range("J34").Value = 4488.1
$ws.Range("K34").Value = 1594.4147
$ws.Range("L34").Value = 4488.1
$ws.Range("M34").Value = -1392.4147
$ws.Range("N34").Value = -4892.1

$ws.Range("H97").Value = 30184.5
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 30184.5
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 30184.5
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -32166.5

$ws.Range("H122").Value = 1815.375
$ws.Range("I122").Value = 1537.1666
$ws.Range("K122").Value = 4611.4998
$ws.Range("M122").Value = -2161.4998

$ws.Range("H132").Value = 1541.4359
$ws.Range("I132").Value = 1337.2285
$ws.Range("J132").Value = 3328.25
$ws.Range("K132").Value = 4011.6855
$ws.Range("L132").Value = 9984.75
$ws.Range("M132").Value = -1481.6855
$ws.Range("N132").Value = -15044.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5452.7144
$ws.Range("I113").Value = 570
$ws.Range("J113").Value = 6601.5884
$ws.Range("K113").Value = 1710
$ws.Range("L113").Value = 19804.7652
$ws.Range("M113").Value = 460
$ws.Range("N113").Value = -24144.7652

$ws.Range("H122").Value = 2658.0962
$ws.Range("I122").Value = 460.48148
$ws.Range("J122").Value = 5031.52
$ws.Range("K122").Value = 4144.33332
$ws.Range("L122").Value = 45283.68000000001
$ws.Range("M122").Value = -1694.33332
$ws.Range("N122").Value = -50183.68000000001

$ws.Range("H131").Value = 10002345
$ws.Range("J131").Value = 11112456
$ws.Range("L131").Value = 33337368
$ws.Range("N131").Value = -33347448

$ws.Range("H133").Value = 3773.2222
$ws.Range("I133").Value = 3249.5
$ws.Range("J133").Value = 3922.8572
$ws.Range("K133").Value = 9748.5
$ws.Range("L133").Value = 11768.5716
$ws.Range("M133").Value = -4688.5
$ws.Range("N133").Value = -21888.5716

$ws.Range("H136").Value = 4504.64
$ws.Range("I136").Value = 1108.3334
$ws.Range("J136").Value = 5577.1577
$ws.Range("K136").Value = 3325.0002
$ws.Range("L136").Value = 16731.4731
$ws.Range("M136").Value = 1774.9998
$ws.Range("N136").Value = -26931.4731

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 53100
$ws.Range("I97").Value = 65125
$ws.Range("K97").Value = 65125
$ws.Range("M97").Value = -64629

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6483.3335
$ws.Range("I100").Value = 7644.4443
$ws.Range("K100").Value = 7644.4443
$ws.Range("M100").Value = -7103.4443

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 101679.5
$ws.Range("I81").Value = 112421.664
$ws.Range("K81").Value = 224843.328
$ws.Range("M81").Value = -223782.328

$ws.Range("H84").Value = 101679.5
$ws.Range("I84").Value = 112421.664
$ws.Range("K84").Value = 1124216.64
$ws.Range("M84").Value = -1118912.64
